$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update each changed cell to match the new crypto price/volume snapshot.
# Numeric-looking values in column D must be force-written as text
# (NumberFormat "@") so Excel keeps them as strings like the original
# inlineStr cells, instead of auto-converting them to real numbers.
# The style is reset back to Normal afterward on each such cell so no
# extra number-format style lingers on the cell.

$ws.Range("D2").Value = "41.525.60"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "2.483.09"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("E7").Value = "  -1.45%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  +2.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0786"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.110"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.84%  "
$ws.Range("D13").Value = "2.867.59"
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.14"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +10.23%  "
$ws.Range("D16").Value = "2.487.74"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.760"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.58%  "
$ws.Range("D18").Value = "41.510.08"
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.49%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0936"
$ws.Range("E20").Value = "  +1.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.18%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("B26").Value = "ImmutableX"
$ws.Range("C26").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.49%  "
$ws.Range("E28").Value = "  -0.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.67"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.47"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.46%  "
$ws.Range("E33").Value = "  +0.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0752"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.52"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.98%  "
$ws.Range("E36").Value = "  -6.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.93"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.105"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.54%  "
$ws.Range("E40").Value = "  +0.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.10"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.18%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.93%  "
$ws.Range("D44").Value = "1.975.48"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0284"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.95"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.34%  "
$ws.Range("D48").Value = "2.725.31"
$ws.Range("E48").Value = "  +1.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.32%  "
